$d = $word.ActiveDocument
$shape = $d.InlineShapes.Item(1)
Write-Host "Height: $($shape.Height)"
Write-Host "Width: $($shape.Width)"
$shape | Get-Member
